# Add 2022-Q3 data
# -----------------------------------------------------------------------
# 1) Build the new "2022-Q3" fund-detail sheet by cloning the existing
#    "2022-Q2" sheet (same header row / styles) right before itself, then
#    trim it down to the 2 fund rows we need and overwrite the content.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)

# the freshly inserted clone lands immediately before the original "2022-Q2"
# sheet, i.e. right after "总计" -- grab it positionally (it is auto-named
# "2022-Q2 (2)", so a name lookup for "2022-Q2" would still hit the original).
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# the clone currently has the same 10 data rows as 2022-Q2 (rows 2-10);
# we only need 2 data rows, so drop the extra ones.
$q3.Rows("4:10").Delete()

# Row 2: 003132 / 德邦新回报灵活配置混合
$q3.Cells.Item(2,1).Value = 0
$q3.Cells.Item(2,2).Value = "'003132"
$q3.Cells.Item(2,3).Value = "德邦新回报灵活配置混合"
$q3.Cells.Item(2,4).Value = "'0.62"
$q3.Cells.Item(2,5).Value = "'72.03"
$q3.Cells.Item(2,6).Value = "'2.35"
$q3.Cells.Item(2,7).Value = "'0.0146"
$q3.Cells.Item(2,8).Value = 9

# Row 3: 080015 / 长盛中小盘精选混合
$q3.Cells.Item(3,1).Value = 1
$q3.Cells.Item(3,2).Value = "'080015"
$q3.Cells.Item(3,3).Value = "长盛中小盘精选混合"
$q3.Cells.Item(3,4).Value = "'0.13"
$q3.Cells.Item(3,5).Value = "'84.41"
$q3.Cells.Item(3,6).Value = "'2.19"
$q3.Cells.Item(3,7).Value = "'0.0028"
$q3.Cells.Item(3,8).Value = 10

# -----------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new leading data row for
#    2022-Q3 and push the existing quarters down by one row.
# -----------------------------------------------------------------------

$total = $wb.Worksheets.Item("总计")

# carry the index-column style down onto the newly used row 5
$total.Cells.Item(2,1).Copy($total.Cells.Item(5,1))

# row 5 (was row 4): 2021-Q4
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q4"
$total.Cells.Item(5,3).Value = 1
$total.Cells.Item(5,4).Value = 0.9

# row 4 (was row 3): 2022-Q1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q1"
$total.Cells.Item(4,3).Value = 11
$total.Cells.Item(4,4).Value = 1.11

# row 3 (was row 2): 2022-Q2
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q2"
$total.Cells.Item(3,3).Value = 9
$total.Cells.Item(3,4).Value = 1.28

# row 2 (new): 2022-Q3
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 2
$total.Cells.Item(2,4).Value = 0.02
